$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 4282.2812
$ws.Cells.Item(40, 9).Value = 3143.92
$ws.Cells.Item(40, 10).Value = 8347.857
$ws.Cells.Item(40, 11).Value = 3143.92
$ws.Cells.Item(40, 12).Value = 8347.857
$ws.Cells.Item(40, 13).Value = -2968.92
$ws.Cells.Item(40, 14).Value = -8697.857
# Row 58
$ws.Cells.Item(58, 8).Value = 1249.1875
$ws.Cells.Item(58, 9).Value = 462.45456
$ws.Cells.Item(58, 10).Value = 2980
$ws.Cells.Item(58, 11).Value = 1387.36368
$ws.Cells.Item(58, 12).Value = 8940
$ws.Cells.Item(58, 13).Value = -1237.36368
$ws.Cells.Item(58, 14).Value = -9240
# Row 137
$ws.Cells.Item(137, 8).Value = 1980264.6
$ws.Cells.Item(137, 9).Value = 3926049
$ws.Cells.Item(137, 10).Value = 1061422
$ws.Cells.Item(137, 11).Value = 11778147
$ws.Cells.Item(137, 12).Value = 3184266
$ws.Cells.Item(137, 13).Value = -11775597
$ws.Cells.Item(137, 14).Value = -3189366

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 632.6667
$ws.Cells.Item(2, 9).Value = 644.7273
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 644.7273
$ws.Cells.Item(2, 12).Value = 500
$ws.Cells.Item(2, 13).Value = -531.7273
$ws.Cells.Item(2, 14).Value = -726
# Row 32
$ws.Cells.Item(32, 8).Value = 1099809.1
$ws.Cells.Item(32, 9).Value = 1114459.9
$ws.Cells.Item(32, 10).Value = 1000
$ws.Cells.Item(32, 11).Value = 1114459.9
$ws.Cells.Item(32, 12).Value = 1000
$ws.Cells.Item(32, 13).Value = -1114172.9
$ws.Cells.Item(32, 14).Value = -1574
# Row 45
$ws.Cells.Item(45, 8).Value = 1339.1111
$ws.Cells.Item(45, 9).Value = 1084.8
$ws.Cells.Item(45, 10).Value = 1657
$ws.Cells.Item(45, 11).Value = 1084.8
$ws.Cells.Item(45, 12).Value = 1657
$ws.Cells.Item(45, 13).Value = -707.8
$ws.Cells.Item(45, 14).Value = -2411
# Row 61
$ws.Cells.Item(61, 8).Value = 348137.47
$ws.Cells.Item(61, 9).Value = 257936.05
$ws.Cells.Item(61, 10).Value = 533287.8
$ws.Cells.Item(61, 11).Value = 257936.05
$ws.Cells.Item(61, 12).Value = 533287.8
$ws.Cells.Item(61, 13).Value = -257724.05
$ws.Cells.Item(61, 14).Value = -533711.8
# Row 74
$ws.Cells.Item(74, 8).Value = 10269060
$ws.Cells.Item(74, 9).Value = 7305550.5
$ws.Cells.Item(74, 11).Value = 7305550.5
$ws.Cells.Item(74, 13).Value = -7304676.5
# Row 77
$ws.Cells.Item(77, 8).Value = 10269060
$ws.Cells.Item(77, 9).Value = 7305550.5
$ws.Cells.Item(77, 11).Value = 36527752.5
$ws.Cells.Item(77, 13).Value = -36523384.5
# Row 116
$ws.Cells.Item(116, 8).Value = 632.6667
$ws.Cells.Item(116, 9).Value = 644.7273
$ws.Cells.Item(116, 10).Value = 500
$ws.Cells.Item(116, 11).Value = 644.7273
$ws.Cells.Item(116, 12).Value = 500
$ws.Cells.Item(116, 13).Value = 1649.2727
$ws.Cells.Item(116, 14).Value = -5088
# Row 122
$ws.Cells.Item(122, 8).Value = 1438
$ws.Cells.Item(122, 9).Value = 996.6667
$ws.Cells.Item(122, 10).Value = 2100
$ws.Cells.Item(122, 11).Value = 2990.0001
$ws.Cells.Item(122, 12).Value = 6300
$ws.Cells.Item(122, 13).Value = -540.0001000000002
$ws.Cells.Item(122, 14).Value = -11200
# Row 136
$ws.Cells.Item(136, 8).Value = 348137.47
$ws.Cells.Item(136, 9).Value = 257936.05
$ws.Cells.Item(136, 10).Value = 533287.8
$ws.Cells.Item(136, 11).Value = 773808.1499999999
$ws.Cells.Item(136, 12).Value = 1599863.4
$ws.Cells.Item(136, 13).Value = -771258.1499999999
$ws.Cells.Item(136, 14).Value = -1604963.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 632.6667
$ws.Cells.Item(3, 9).Value = 644.7273
$ws.Cells.Item(3, 10).Value = 500
$ws.Cells.Item(3, 11).Value = 644.7273
$ws.Cells.Item(3, 12).Value = 500
$ws.Cells.Item(3, 13).Value = -530.7273
$ws.Cells.Item(3, 14).Value = -728
# Row 80
$ws.Cells.Item(80, 8).Value = 689
$ws.Cells.Item(80, 10).Value = 525
$ws.Cells.Item(80, 12).Value = 525
$ws.Cells.Item(80, 14).Value = -2521
# Row 83
$ws.Cells.Item(83, 8).Value = 689
$ws.Cells.Item(83, 10).Value = 525
$ws.Cells.Item(83, 12).Value = 2625
$ws.Cells.Item(83, 14).Value = -12609
# Row 86
$ws.Cells.Item(86, 8).Value = 4541.6523
$ws.Cells.Item(86, 9).Value = 5034.263
$ws.Cells.Item(86, 11).Value = 5034.263
$ws.Cells.Item(86, 13).Value = -3911.263
# Row 89
$ws.Cells.Item(89, 8).Value = 4541.6523
$ws.Cells.Item(89, 9).Value = 5034.263
$ws.Cells.Item(89, 11).Value = 25171.315
$ws.Cells.Item(89, 13).Value = -19555.315

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 934.3077
$ws.Cells.Item(16, 9).Value = 834.5
$ws.Cells.Item(16, 10).Value = 1019.8571
$ws.Cells.Item(16, 11).Value = 834.5
$ws.Cells.Item(16, 12).Value = 1019.8571
$ws.Cells.Item(16, 13).Value = -547.5
$ws.Cells.Item(16, 14).Value = -1593.8571
# Row 31
$ws.Cells.Item(31, 8).Value = 1512.77
$ws.Cells.Item(31, 9).Value = 669.1724
$ws.Cells.Item(31, 10).Value = 2677.738
$ws.Cells.Item(31, 11).Value = 669.1724
$ws.Cells.Item(31, 12).Value = 2677.738
$ws.Cells.Item(31, 13).Value = -374.1724
$ws.Cells.Item(31, 14).Value = -3267.738
# Row 34
$ws.Cells.Item(34, 8).Value = 1512.77
$ws.Cells.Item(34, 9).Value = 669.1724
$ws.Cells.Item(34, 10).Value = 2677.738
$ws.Cells.Item(34, 11).Value = 669.1724
$ws.Cells.Item(34, 12).Value = 2677.738
$ws.Cells.Item(34, 13).Value = -467.1724
$ws.Cells.Item(34, 14).Value = -3081.738
# Row 58
$ws.Cells.Item(58, 8).Value = 3472.3508
$ws.Cells.Item(58, 9).Value = 4271.793
$ws.Cells.Item(58, 10).Value = 2644.3572
$ws.Cells.Item(58, 11).Value = 4271.793
$ws.Cells.Item(58, 12).Value = 2644.3572
$ws.Cells.Item(58, 13).Value = -4068.793
$ws.Cells.Item(58, 14).Value = -3050.3572
# Row 113
$ws.Cells.Item(113, 8).Value = 934.3077
$ws.Cells.Item(113, 9).Value = 834.5
$ws.Cells.Item(113, 10).Value = 1019.8571
$ws.Cells.Item(113, 11).Value = 834.5
$ws.Cells.Item(113, 12).Value = 1019.8571
$ws.Cells.Item(113, 13).Value = 1335.5
$ws.Cells.Item(113, 14).Value = -5359.8571
# Row 122
$ws.Cells.Item(122, 8).Value = 1856.2
$ws.Cells.Item(122, 9).Value = 2434.1667
$ws.Cells.Item(122, 10).Value = 989.25
$ws.Cells.Item(122, 11).Value = 7302.500100000001
$ws.Cells.Item(122, 12).Value = 2967.75
$ws.Cells.Item(122, 13).Value = -4852.500100000001
$ws.Cells.Item(122, 14).Value = -7867.75
# Row 134
$ws.Cells.Item(134, 8).Value = 11629102
$ws.Cells.Item(134, 9).Value = 15152364
$ws.Cells.Item(134, 10).Value = 2340
$ws.Cells.Item(134, 11).Value = 45457092
$ws.Cells.Item(134, 12).Value = 7020
$ws.Cells.Item(134, 13).Value = -45454557
$ws.Cells.Item(134, 14).Value = -12090
# Row 136
$ws.Cells.Item(136, 8).Value = 3472.3508
$ws.Cells.Item(136, 9).Value = 4271.793
$ws.Cells.Item(136, 10).Value = 2644.3572
$ws.Cells.Item(136, 11).Value = 12815.379
$ws.Cells.Item(136, 12).Value = 7933.071599999999
$ws.Cells.Item(136, 13).Value = -10265.379
$ws.Cells.Item(136, 14).Value = -13033.0716

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1014.80554
$ws.Cells.Item(5, 9).Value = 776.46155
$ws.Cells.Item(5, 10).Value = 1149.5217
$ws.Cells.Item(5, 11).Value = 2329.38465
$ws.Cells.Item(5, 12).Value = 3448.5651
$ws.Cells.Item(5, 13).Value = -2217.38465
$ws.Cells.Item(5, 14).Value = -3672.5651
# Row 33
$ws.Cells.Item(33, 8).Value = 293.14285
$ws.Cells.Item(33, 9).Value = 211.6923
$ws.Cells.Item(33, 10).Value = 425.5
$ws.Cells.Item(33, 11).Value = 1270.1538
$ws.Cells.Item(33, 12).Value = 2553
$ws.Cells.Item(33, 13).Value = -987.1538
$ws.Cells.Item(33, 14).Value = -3119
# Row 107
$ws.Cells.Item(107, 8).Value = 1174.5588
$ws.Cells.Item(107, 9).Value = 482.84616
$ws.Cells.Item(107, 10).Value = 1602.762
$ws.Cells.Item(107, 11).Value = 1448.53848
$ws.Cells.Item(107, 12).Value = 4808.286
$ws.Cells.Item(107, 13).Value = 471.4615200000001
$ws.Cells.Item(107, 14).Value = -8648.286
# Row 131
$ws.Cells.Item(131, 8).Value = 1188.2969
$ws.Cells.Item(131, 9).Value = 1258
$ws.Cells.Item(131, 10).Value = 1175.3889
$ws.Cells.Item(131, 11).Value = 3774
$ws.Cells.Item(131, 12).Value = 3526.1667
$ws.Cells.Item(131, 13).Value = 1266
$ws.Cells.Item(131, 14).Value = -13606.1667
# Row 132
$ws.Cells.Item(132, 8).Value = 1747.2424
$ws.Cells.Item(132, 9).Value = 3843.4285
$ws.Cells.Item(132, 10).Value = 1182.8846
$ws.Cells.Item(132, 11).Value = 34590.8565
$ws.Cells.Item(132, 12).Value = 10645.9614
$ws.Cells.Item(132, 13).Value = -32060.8565
$ws.Cells.Item(132, 14).Value = -15705.9614
# Row 133
$ws.Cells.Item(133, 8).Value = 4628.5713
$ws.Cells.Item(133, 9).Value = 5010
$ws.Cells.Item(133, 10).Value = 4342.5
$ws.Cells.Item(133, 11).Value = 15030
$ws.Cells.Item(133, 12).Value = 13027.5
$ws.Cells.Item(133, 13).Value = -9970
$ws.Cells.Item(133, 14).Value = -23147.5
# Row 135
$ws.Cells.Item(135, 8).Value = 1014.80554
$ws.Cells.Item(135, 9).Value = 776.46155
$ws.Cells.Item(135, 10).Value = 1149.5217
$ws.Cells.Item(135, 11).Value = 6988.15395
$ws.Cells.Item(135, 12).Value = 10345.6953
$ws.Cells.Item(135, 13).Value = -4453.15395
$ws.Cells.Item(135, 14).Value = -15415.6953

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Cells.Item(113, 8).Value = 35722660
$ws.Cells.Item(113, 9).Value = 71443704
$ws.Cells.Item(113, 10).Value = 1614.2858
$ws.Cells.Item(113, 11).Value = 71443704
$ws.Cells.Item(113, 12).Value = 1614.2858
$ws.Cells.Item(113, 13).Value = -71441534
$ws.Cells.Item(113, 14).Value = -5954.2858
# Row 122
$ws.Cells.Item(122, 8).Value = 1222.0834
$ws.Cells.Item(122, 9).Value = 1038
$ws.Cells.Item(122, 10).Value = 1353.5714
$ws.Cells.Item(122, 11).Value = 3114
$ws.Cells.Item(122, 12).Value = 4060.7142
$ws.Cells.Item(122, 13).Value = -664
$ws.Cells.Item(122, 14).Value = -8960.7142
# Row 126
$ws.Cells.Item(126, 8).Value = 3212.5
$ws.Cells.Item(126, 9).Value = 2900
$ws.Cells.Item(126, 10).Value = 3400
$ws.Cells.Item(126, 11).Value = 8700
$ws.Cells.Item(126, 12).Value = 10200
$ws.Cells.Item(126, 13).Value = -6230
$ws.Cells.Item(126, 14).Value = -15140
# Row 132
$ws.Cells.Item(132, 8).Value = 2720790.2
$ws.Cells.Item(132, 9).Value = 4169795
$ws.Cells.Item(132, 10).Value = 3906.125
$ws.Cells.Item(132, 11).Value = 12509385
$ws.Cells.Item(132, 12).Value = 11718.375
$ws.Cells.Item(132, 13).Value = -12506855
$ws.Cells.Item(132, 14).Value = -16778.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 1264.5581
$ws.Cells.Item(46, 9).Value = 1094.1724
$ws.Cells.Item(46, 11).Value = 1094.1724
$ws.Cells.Item(46, 13).Value = -906.1723999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 17558012
$ws.Cells.Item(132, 9).Value = 7083.1665
$ws.Cells.Item(132, 10).Value = 25658442
$ws.Cells.Item(132, 11).Value = 21249.4995
$ws.Cells.Item(132, 12).Value = 76975326
$ws.Cells.Item(132, 13).Value = -18719.4995
$ws.Cells.Item(132, 14).Value = -76980386
# Row 136
$ws.Cells.Item(136, 8).Value = 3003.125
$ws.Cells.Item(136, 9).Value = 1508.55
$ws.Cells.Item(136, 10).Value = 6739.5625
$ws.Cells.Item(136, 11).Value = 4525.65
$ws.Cells.Item(136, 12).Value = 20218.6875
$ws.Cells.Item(136, 13).Value = -1975.65
$ws.Cells.Item(136, 14).Value = -25318.6875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 646.6667
$ws.Cells.Item(107, 9).Value = 470
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 11).Value = 1410
$ws.Cells.Item(107, 12).Value = 3000
$ws.Cells.Item(107, 13).Value = 510
$ws.Cells.Item(107, 14).Value = -6840
# Row 113
$ws.Cells.Item(113, 8).Value = 330.10526
$ws.Cells.Item(113, 9).Value = 333.64706
$ws.Cells.Item(113, 10).Value = 300
$ws.Cells.Item(113, 11).Value = 1000.94118
$ws.Cells.Item(113, 12).Value = 900
$ws.Cells.Item(113, 13).Value = 1169.05882
$ws.Cells.Item(113, 14).Value = -5240
# Row 136
$ws.Cells.Item(136, 8).Value = 20401902
$ws.Cells.Item(136, 9).Value = 24660192
$ws.Cells.Item(136, 10).Value = 7931193.5
$ws.Cells.Item(136, 11).Value = 73980576
$ws.Cells.Item(136, 12).Value = 23793580.5
$ws.Cells.Item(136, 13).Value = -73978026
$ws.Cells.Item(136, 14).Value = -23798680.5
